# Final commit of upload excel file
# - Corrects a handful of contact values in the "contacts" sheet
# - Tidies the "Hobbies" entries (drop trailing comma)
# - Bumps the row height of the three data rows slightly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (first contact) corrections
$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"

# Row 3 (second contact) corrections
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# Slightly taller header/data rows
$ws.Rows("1:3").RowHeight = 19.5
